$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("branchAndBound")

# --- Update the workbook window view ---
$wb.Windows.Item(1).Left = 39660
$wb.Windows.Item(1).Top = 2780
$wb.Windows.Item(1).Width = 27180
$wb.Windows.Item(1).Height = 14440

# --- Update the selected cell on the branchAndBound sheet ---
[void]$ws.Range("G14").Select()

# --- Row 10 updates ---
$ws.Range("D10").Value = -1340.51
$ws.Range("G10").Value = 8103
$ws.Range("H10").Value = 1711
$ws.Range("I10").Value = 381

# --- Row 11 updates ---
$ws.Range("D11").Value = -1350.16
$ws.Range("G11").Value = 9317
$ws.Range("H11").Value = 809
$ws.Range("I11").Value = 39

# --- Row 12 updates ---
$ws.Range("D12").Value = -1356.72
$ws.Range("G12").Value = 10220
$ws.Range("H12").Value = 1305
$ws.Range("I12").Value = 7

# --- Row 18: fill in previously-empty result cells ---
$ws.Range("C18").Value = -1338.68
$ws.Range("D18").Value = -1341.93
$ws.Range("E18").Formula = "=(D18-C18)/D18"
$ws.Range("F18").Value = ">20,000"
$ws.Range("G18").Value = 9924
$ws.Range("I18").Value = 363

# --- Row 19: fill in previously-empty result cells ---
$ws.Range("C19").Value = -1313.93
$ws.Range("D19").Value = -1414.68
$ws.Range("E19").Formula = "=(D19-C19)/D19"
$ws.Range("F19").Value = ">20,000"
$ws.Range("G19").Value = 10511
$ws.Range("I19").Value = 43

# --- Row 20: fill in previously-empty result cells ---
$ws.Range("C20").Value = -1308.9100000000001
$ws.Range("D20").Value = -1491.52
$ws.Range("E20").Formula = "=(D20-C20)/D20"
$ws.Range("F20").Value = ">20,000"
$ws.Range("G20").Value = 10893
$ws.Range("I20").Value = 9

# --- Now change the shared text from ">20,000" to ">10,000" everywhere it's used ---
$ws.Range("F10").Value = ">10,000"
$ws.Range("F11").Value = ">10,000"
$ws.Range("F12").Value = ">10,000"
$ws.Range("F18").Value = ">10,000"
$ws.Range("F19").Value = ">10,000"
$ws.Range("F20").Value = ">10,000"
